$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Investor 1"
$ws.Range("A3").Value = "Investor 2"
$ws.Range("A4").Value = "Investor 3"
$ws.Range("A5").Value = "Investor 4"
$ws.Range("A6").Value = "Investor 5"

$ws.Range("A2:A6").Select()
